$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

function Set-TextValue($cell, $value, $styleTemplateCell) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $styleTemplateCell.Style
}

# Row 9: date / who / defect reason / count (mirrors rows 4-8 above it)
Set-TextValue $ws.Cells.Item(9,1) "30/01/2018" $ws.Cells.Item(5,1)
Set-TextValue $ws.Cells.Item(9,2) "3012" $ws.Cells.Item(5,2)
Set-TextValue $ws.Cells.Item(9,3) "Не вірна довжина проводу" $ws.Cells.Item(7,3)
Set-TextValue $ws.Cells.Item(9,4) "1" $ws.Cells.Item(5,4)

# Row 10
Set-TextValue $ws.Cells.Item(10,1) "31/01/2018" $ws.Cells.Item(5,1)
Set-TextValue $ws.Cells.Item(10,2) "1" $ws.Cells.Item(5,2)
Set-TextValue $ws.Cells.Item(10,3) "Пошкодження поверхні контакту" $ws.Cells.Item(7,3)
Set-TextValue $ws.Cells.Item(10,4) "1" $ws.Cells.Item(5,4)

# Row 11
Set-TextValue $ws.Cells.Item(11,1) "**" $ws.Cells.Item(5,1)

# Remove the empty placeholder rows 25-27 (they hold no data; row 28 keeps its row number)
$ws.Rows(25).ClearContents()
$ws.Rows(26).ClearContents()
$ws.Rows(27).ClearContents()
